# guide41_dashboad.xlsx edit:
#  - Split the combined "p2" page (language-menu + preferred-language
#    instructions) into two pages:
#      * p2 ("常時英語でmoodleを利用したい")  -> keeps the preferred-language section
#      * p3 ("一時的に英語でmoodleを利用したい") -> new page with the language-menu section
#
# Strategy: duplicate the existing p2 sheet to create p3, then trim each
# sheet down to the rows it should keep and fix up the text that changed.

$wb = $excel.ActiveWorkbook

$idx = $wb.Worksheets.Item("index")
$p1  = $wb.Worksheets.Item("p1")
$p2  = $wb.Worksheets.Item("p2")

# ---------------------------------------------------------------------
# 1. Create "p3" as an exact copy of "p2" (so it inherits every style,
#    column width, row height, etc. that the new page needs), placed
#    right after "p2".
# ---------------------------------------------------------------------
$p2.Copy($null, $p2) | Out-Null
$p3 = $wb.Worksheets.Item($p2.Index + 1)
$p3.Name = "p3"

# ---------------------------------------------------------------------
# 2. Trim "p3" down to the "language menu" section (old rows 1-9) and
#    update its text.
# ---------------------------------------------------------------------
$lastRow3 = $p3.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11
if ($lastRow3 -gt 9) {
    $p3.Rows("10:" + $lastRow3).Delete() | Out-Null
}

$p3.Range("B2").Value = "一時的に英語でmoodleを利用したい"
$p3.Range("B9").Value = "表示したい言語を選択します`nデフォルトでは、日本語（ja）になっているので、English(en)を選択します"
$p3.Rows.Item(9).RowHeight = 27

$p3.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Trim "p2" down to the "preferred language" section by removing the
#    old "language menu" rows (7-9); everything below shifts up and
#    keeps its original formatting.
# ---------------------------------------------------------------------
$p2.Rows("7:9").Delete() | Out-Null

$p2.Range("B2").Value = "常時英語でmoodleを利用したい"
$p2.Range("B10").Value = "優先言語の▼をクリックし、リストを表示させ優先させたい言語を選択します`nデフォルトでは、日本語（ja）が優先言語になっているので、English(en)を選択します"
$p2.Range("B11").Value = "［変更を保存する］をクリックします"
$p2.Rows.Item(10).RowHeight = 27

$p2.Range("B15").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Misc cursor/selection bookkeeping on the other sheets and make the
#    new page the active tab, matching the saved workbook state.
# ---------------------------------------------------------------------
$p1.Activate() | Out-Null
$p1.Range("B7").Select() | Out-Null

$p3.Activate() | Out-Null
